$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update forecast values in columns C (row-1 offset) and E (row offset) per fix to naive component forecaster
$ws.Cells.Item(2, 3).Value = 3.123541145015474
$ws.Cells.Item(2, 5).Value = 3.275761698145385
$ws.Cells.Item(3, 3).Value = 3.959010658874851
$ws.Cells.Item(3, 5).Value = 3.716811705074696
$ws.Cells.Item(4, 3).Value = 4.722695063536686
$ws.Cells.Item(4, 5).Value = 5.522497812297966
$ws.Cells.Item(5, 3).Value = 6.739021039846627
$ws.Cells.Item(5, 5).Value = 3.88772167754905
$ws.Cells.Item(6, 3).Value = 2.619839412265601
$ws.Cells.Item(6, 5).Value = 2.645362710332533
$ws.Cells.Item(7, 3).Value = -0.7919564768266385
$ws.Cells.Item(7, 5).Value = 0.8976572162586516
$ws.Cells.Item(8, 3).Value = 1.877689851450803
$ws.Cells.Item(8, 5).Value = 1.899791808163398
$ws.Cells.Item(9, 3).Value = 2.705004599189187
$ws.Cells.Item(9, 5).Value = 2.111643953433728
$ws.Cells.Item(10, 3).Value = 1.110374544249249
$ws.Cells.Item(10, 5).Value = 2.576021643263426
$ws.Cells.Item(11, 3).Value = 2.267566233338814
$ws.Cells.Item(11, 5).Value = 2.348556921565126
$ws.Cells.Item(12, 3).Value = 2.688433258834588
$ws.Cells.Item(12, 5).Value = 2.730428819177333
$ws.Cells.Item(13, 3).Value = 1.014079695989589
$ws.Cells.Item(13, 5).Value = 1.901826179618205
$ws.Cells.Item(14, 3).Value = 3.013853578092252
$ws.Cells.Item(14, 5).Value = 2.435357021275819
$ws.Cells.Item(15, 3).Value = 1.331295149770684
$ws.Cells.Item(15, 5).Value = 1.127665471558248
$ws.Cells.Item(16, 3).Value = 0.04589006555719699
$ws.Cells.Item(16, 5).Value = 0.8401596151991431
$ws.Cells.Item(17, 3).Value = 0.009546395482029624
$ws.Cells.Item(17, 5).Value = 0.02245646656315881
$ws.Cells.Item(18, 3).Value = 0.8709390141433015
$ws.Cells.Item(18, 5).Value = 0.7407001102931465
$ws.Cells.Item(19, 3).Value = 0.7652063367885598
$ws.Cells.Item(19, 5).Value = 1.397861196490657
$ws.Cells.Item(20, 3).Value = 2.267579219134386
$ws.Cells.Item(20, 5).Value = 2.220975586034668
$ws.Cells.Item(21, 3).Value = 3.146753122914103
$ws.Cells.Item(21, 5).Value = 2.456586080053058
$ws.Cells.Item(22, 3).Value = 1.769033835366818
$ws.Cells.Item(22, 5).Value = 1.086632508372576
$ws.Cells.Item(23, 3).Value = -4.774715709990263
$ws.Cells.Item(23, 5).Value = -0.8529145826070339
$ws.Cells.Item(24, 3).Value = 1.95493704440024
$ws.Cells.Item(24, 5).Value = 1.373265374526711
$ws.Cells.Item(25, 3).Value = 3.478075069442799
$ws.Cells.Item(25, 5).Value = 2.807281147895924
$ws.Cells.Item(26, 3).Value = 1.232342134690434
$ws.Cells.Item(26, 5).Value = 2.058767060726563
$ws.Cells.Item(27, 3).Value = 0.2542811494408159
$ws.Cells.Item(27, 5).Value = 1.100861823237564
$ws.Cells.Item(28, 3).Value = 1.519778766382096
$ws.Cells.Item(28, 5).Value = 0.6176362615319508
$ws.Cells.Item(29, 3).Value = 1.469441753880329
$ws.Cells.Item(29, 5).Value = 1.593307036690472
$ws.Cells.Item(30, 3).Value = 1.638203081492495
$ws.Cells.Item(30, 5).Value = 1.425381384581903
$ws.Cells.Item(31, 3).Value = 2.268697431234346
$ws.Cells.Item(31, 5).Value = 2.525738252590148
$ws.Cells.Item(32, 3).Value = 1.984425467899631
$ws.Cells.Item(32, 5).Value = 1.272991730262341
$ws.Cells.Item(33, 3).Value = 0.6066448776129052
$ws.Cells.Item(33, 5).Value = 0.8323605077169782
$ws.Cells.Item(34, 3).Value = -4.243076347305386
$ws.Cells.Item(34, 5).Value = -0.2172115314521883
$ws.Cells.Item(35, 3).Value = 1.438499295329754
$ws.Cells.Item(35, 5).Value = 1.725839624491665
$ws.Cells.Item(36, 3).Value = 1.906593537051537
$ws.Cells.Item(36, 5).Value = 1.222362529774923
$ws.Cells.Item(37, 3).Value = 0.08348019664223827
$ws.Cells.Item(37, 5).Value = 0.9116426337375527
$ws.Cells.Item(38, 3).Value = -0.214505326882275
$ws.Cells.Item(38, 5).Value = -0.02742999272021818
$ws.Cells.Item(39, 3).Value = 0.1651547428133782
$ws.Cells.Item(39, 5).Value = -0.1331392688890709
